$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-51 refresh: Coin Price / Volume(1h) snapshot update.
# Rows 46/47 additionally swap (EnergySwap <-> Decentraland re-ranked).
$updates = @(
    @{ Addr = "D2"; Value = "28.964.61" },
    @{ Addr = "E2"; Value = "  +5.45%  " },
    @{ Addr = "D3"; Value = "1.915.76" },
    @{ Addr = "E3"; Value = "  +4.61%  " },
    @{ Addr = "D4"; Value = "1.001" },
    @{ Addr = "E4"; Value = "  -0.14%  " },
    @{ Addr = "D5"; Value = "338.45" },
    @{ Addr = "E5"; Value = "  +2.00%  " },
    @{ Addr = "E6"; Value = "  -0.17%  " },
    @{ Addr = "D7"; Value = "0.4741" },
    @{ Addr = "E7"; Value = "  +3.52%  " },
    @{ Addr = "D8"; Value = "0.4062" },
    @{ Addr = "E8"; Value = "  +6.71%  " },
    @{ Addr = "D9"; Value = "48.12" },
    @{ Addr = "E9"; Value = "  +3.41%  " },
    @{ Addr = "D10"; Value = "0.08173" },
    @{ Addr = "E10"; Value = "  +3.53%  " },
    @{ Addr = "E11"; Value = "  +6.41%  " },
    @{ Addr = "D12"; Value = "22.45" },
    @{ Addr = "E12"; Value = "  +6.55%  " },
    @{ Addr = "D13"; Value = "1.897.47" },
    @{ Addr = "E13"; Value = "  +3.99%  " },
    @{ Addr = "D14"; Value = "6.084" },
    @{ Addr = "E14"; Value = "  +3.34%  " },
    @{ Addr = "D15"; Value = "7.387" },
    @{ Addr = "E15"; Value = "  +4.44%  " },
    @{ Addr = "D16"; Value = "91.75" },
    @{ Addr = "E16"; Value = "  +2.75%  " },
    @{ Addr = "D17"; Value = "1.002" },
    @{ Addr = "E17"; Value = "  -0.11%  " },
    @{ Addr = "E18"; Value = "  +2.81%  " },
    @{ Addr = "D19"; Value = "0.06633" },
    @{ Addr = "E19"; Value = "  +0.34%  " },
    @{ Addr = "D20"; Value = "17.88" },
    @{ Addr = "E20"; Value = "  +4.33%  " },
    @{ Addr = "D21"; Value = "1.000" },
    @{ Addr = "E21"; Value = "  -0.21%  " },
    @{ Addr = "D22"; Value = "29.003.84" },
    @{ Addr = "E22"; Value = "  +5.64%  " },
    @{ Addr = "D23"; Value = "5.588" },
    @{ Addr = "E23"; Value = "  +4.58%  " },
    @{ Addr = "D24"; Value = "11.18" },
    @{ Addr = "E24"; Value = "  +3.43%  " },
    @{ Addr = "E25"; Value = "  -1.03%  " },
    @{ Addr = "D26"; Value = "2.138.39" },
    @{ Addr = "E26"; Value = "  +4.76%  " },
    @{ Addr = "D27"; Value = "160.83" },
    @{ Addr = "D28"; Value = "20.06" },
    @{ Addr = "E28"; Value = "  +3.41%  " },
    @{ Addr = "D29"; Value = "2.188" },
    @{ Addr = "E29"; Value = "  +5.80%  " },
    @{ Addr = "D30"; Value = "5.536" },
    @{ Addr = "E30"; Value = "  +4.53%  " },
    @{ Addr = "D31"; Value = "121.19" },
    @{ Addr = "E31"; Value = "  +2.12%  " },
    @{ Addr = "D32"; Value = "1.015" },
    @{ Addr = "E32"; Value = "  +7.73%  " },
    @{ Addr = "D33"; Value = "0.09576" },
    @{ Addr = "E33"; Value = "  +2.88%  " },
    @{ Addr = "D34"; Value = "1.440" },
    @{ Addr = "E34"; Value = "  +8.09%  " },
    @{ Addr = "D35"; Value = "3.639" },
    @{ Addr = "E35"; Value = "  +1.71%  " },
    @{ Addr = "D36"; Value = "5.433" },
    @{ Addr = "E36"; Value = "  +3.63%  " },
    @{ Addr = "D37"; Value = "0.06212" },
    @{ Addr = "E37"; Value = "  +4.84%  " },
    @{ Addr = "D38"; Value = "8.722" },
    @{ Addr = "E38"; Value = "  +7.69%  " },
    @{ Addr = "D39"; Value = "0.02287" },
    @{ Addr = "E39"; Value = "  +4.62%  " },
    @{ Addr = "D40"; Value = "1.205" },
    @{ Addr = "E40"; Value = "  +4.75%  " },
    @{ Addr = "D41"; Value = "0.6039" },
    @{ Addr = "E41"; Value = "  +4.50%  " },
    @{ Addr = "D42"; Value = "10.58" },
    @{ Addr = "E42"; Value = "  +5.69%  " },
    @{ Addr = "E43"; Value = "  +4.05%  " },
    @{ Addr = "E44"; Value = "  -0.11%  " },
    @{ Addr = "D45"; Value = "1.283" },
    @{ Addr = "E45"; Value = "  +1.50%  " },
    @{ Addr = "B46"; Value = "EnergySwap" },
    @{ Addr = "C46"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" },
    @{ Addr = "D46"; Value = "12.40" },
    @{ Addr = "E46"; Value = "  +3.84%  " },
    @{ Addr = "B47"; Value = "Decentraland" },
    @{ Addr = "C47"; Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana" },
    @{ Addr = "D47"; Value = "0.5628" },
    @{ Addr = "E47"; Value = "  +3.10%  " },
    @{ Addr = "D48"; Value = "1.982" },
    @{ Addr = "E48"; Value = "  +5.95%  " },
    @{ Addr = "D49"; Value = "0.07243" },
    @{ Addr = "E49"; Value = "  +9.59%  " },
    @{ Addr = "D50"; Value = "2.150" },
    @{ Addr = "E50"; Value = "  +18.50%  " },
    @{ Addr = "D51"; Value = "113.21" },
    @{ Addr = "E51"; Value = "  +2.11%  " }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    $val = $u.Value
    # Price column (D) values are plain decimal-looking strings (e.g. "1.001",
    # "338.45") that must stay text, not be coerced into numbers on entry -
    # force text format first so they round-trip exactly like the source data.
    if ($u.Addr -match '^D' -and $val -match '^-?[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}
